# Daily attendance processing - 2025-11-02 13:06:39
# Reorders the "Recorded By" names in column G so that system-type
# recorders ("System" / "system") are listed before the human recorder
# email address, for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $value = $cell.Value2

    if ($value -eq $null) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value.IndexOf(",") -lt 0) { continue }

    $parts = $value -split ", "

    if ($parts.Count -eq 2 -and $parts[0] -eq "dnasr281@gmail.com" -and $parts[1] -eq "System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($parts.Count -eq 3 -and $parts[0] -eq "backup@backdoor.com" -and $parts[1] -eq "system" -and $parts[2] -eq "System") {
        $cell.Value2 = "system, backup@backdoor.com, System"
    }
}
